$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the new columns.
#    Before:  A..H | I=Turnaround J=Throughput K=CPU L=Fairness M=CtxSwitch
#    Insert a blank column at I  -> shifts Turnaround..CtxSwitch to J..N
#    Insert two blank columns at N -> shifts CtxSwitch (now at N) to P
# ---------------------------------------------------------------------
$ws.Columns("I").Insert()
$ws.Range("N1:O1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. New "Elapsed Time" column (I)
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "Elapsed Time"
$ws.Range("I2:I11").Formula = "=(C2-`$B`$12)/100"

# ---------------------------------------------------------------------
# 3. Rename / re-purpose the fairness columns.
#    M currently still holds the old "Fairness" header/value (STDEV.P of
#    Real_time D2:D11); that becomes "Fairness(RT)" and moves to N.
#    M becomes "Fairness(QWT)" with a new STDEV.P(G2:G11) formula.
#    O becomes the new "Fairness(ET)" = STDEV.P(I2:I11).
# ---------------------------------------------------------------------
$ws.Range("N1").Value = $ws.Range("M1").Text
$ws.Range("N12").Formula = $ws.Range("M12").Formula

$ws.Range("M1").Value = "Fairness(QWT)"
$ws.Range("M12").Formula = "=_xlfn.STDEV.P(G2:G11)"

$ws.Range("O1").Value = "Fairness(ET)"
$ws.Range("O12").Formula = "=_xlfn.STDEV.P(I2:I11)"

# ---------------------------------------------------------------------
# 4. Cosmetics to match the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("M18").Select()
$ws.Rows("1:12").RowHeight = 14.25

$wb.Theme.ThemeFontScheme.MinorFont.Item(1).Name = "Arial"
$wb.Theme.ThemeFontScheme.MajorFont.Item(1).Name = "Times New Roman"

$excel.Calculate()
